$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from an existing header cell (style s="1") onto the
# new header cell H1, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill in the new "Save" column (H) values for rows 2-22.
$saveValues = @(1,0,1,1,0,1,1,1,0,1,1,1,1,1,0,1,1,0,1,1,1)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
